# Weekly update: insert a new price record as row 11 (pushing the existing
# rows 11-32 down to 12-33), matching the "Fruta / hortaliza, semanal" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 11; this shifts old rows 11..32 down to 12..33
$ws.Rows("11:11").Insert()

# Populate the new row 11 with this week's record
$ws.Cells.Item(11, 1).Value2  = 1
$ws.Cells.Item(11, 2).Value2  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(11, 3).Value2  = "Arica y Parinacota"
$ws.Cells.Item(11, 4).Value2  = 45002
$ws.Cells.Item(11, 5).Value2  = 15
$ws.Cells.Item(11, 6).Value2  = "Fruta"
$ws.Cells.Item(11, 7).Value2  = 100103
$ws.Cells.Item(11, 8).Value2  = "Frutos de hueso (carozo)"
$ws.Cells.Item(11, 9).Value2  = 100103002
$ws.Cells.Item(11, 10).Value2 = "Ciruela"
$ws.Cells.Item(11, 11).Value2 = "Angeleno"
$ws.Cells.Item(11, 12).Value2 = "Segunda"
$ws.Cells.Item(11, 13).Value2 = 300
$ws.Cells.Item(11, 14).Value2 = 21000
$ws.Cells.Item(11, 15).Value2 = 22000
$ws.Cells.Item(11, 16).Value2 = 21500
$ws.Cells.Item(11, 17).Value2 = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(11, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(11, 19).Value2 = 1194
$ws.Cells.Item(11, 20).Value2 = 18
